$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the formula in C4 (now multiplied by 2 twice)
$ws.Range("C4").Formula = "=C3*2*2"

# Add a new row of data: Person 4 / QA / 2
$ws.Range("A5").Value = "Person 4 "
$ws.Range("B5").Value = "QA"
$ws.Range("C5").Value = 2

# Copy style (red font + black fill) from existing row to the new row
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)  # xlPasteFormats

# Update the selected cell in the sheet view
$ws.Range("D5").Select()
